$d = $word.ActiveDocument

# --- Collapse split runs in Title / Author / Abstract paragraphs into single runs ---
$xTitle = $d.Paragraphs.Item(1).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:pStyle w:val=''Title''/></w:pPr><w:r><w:t xml:space=''preserve''>Questions: Laws of indices</w:t></w:r></w:p>')
$xAuthor = $d.Paragraphs.Item(2).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:pStyle w:val=''Author''/></w:pPr><w:r><w:t xml:space=''preserve''>Isabella Lewis, Akshat Srivastava</w:t></w:r></w:p>')
$xAbstract = $d.Paragraphs.Item(4).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:pStyle w:val=''Abstract''/></w:pPr><w:r><w:t xml:space=''preserve''>A selection of questions for the study guide on laws of indices.</w:t></w:r></w:p>')

# --- Reorder m:dPr children (sepChr before endChr) inside the bracket/delimiter math objects ---
$xOMath10 = $d.OMaths.Item(10).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>5</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath11 = $d.OMaths.Item(11).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>35</m:t></m:r></m:e><m:sup><m:r><m:t>0</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:e><m:sup><m:r><m:t>9</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath12 = $d.OMaths.Item(12).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>35</m:t></m:r></m:e><m:sup><m:r><m:t>9</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:e><m:sup><m:r><m:t>0</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath13 = $d.OMaths.Item(13).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>729</m:t></m:r></m:e><m:sup><m:r><m:t>9</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:e><m:sup><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>9</m:t></m:r></m:den></m:f></m:sup></m:sSup></m:oMath>')
$xOMath18 = $d.OMaths.Item(18).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:sSup><m:e><m:r><m:t>4</m:t></m:r></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:r><m:t>3</m:t></m:r></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup></m:num><m:den><m:sSup><m:e><m:r><m:t>6</m:t></m:r></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup></m:den></m:f></m:e></m:d></m:oMath>')
$xOMath19 = $d.OMaths.Item(19).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:sSup><m:e><m:r><m:t>4</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:r><m:t>8</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:num><m:den><m:sSup><m:e><m:r><m:t>2</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:den></m:f></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath20 = $d.OMaths.Item(20).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:d><m:dPr><m:begChr m:val="[" /><m:sepChr m:val="" /><m:endChr m:val="]" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>5</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:num><m:den><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>2</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r></m:sup></m:sSup></m:den></m:f></m:oMath>')
$xOMath21 = $d.OMaths.Item(21).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>3</m:t></m:r></m:num><m:den><m:r><m:t>5</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup></m:num><m:den><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>8</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup></m:den></m:f></m:oMath>')
$xOMath22 = $d.OMaths.Item(22).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>2</m:t></m:r></m:num><m:den><m:r><m:t>3</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>14</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>9</m:t></m:r></m:num><m:den><m:r><m:t>12</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>14</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath23 = $d.OMaths.Item(23).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>b</m:t></m:r></m:e><m:sup><m:r><m:t>7</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath27 = $d.OMaths.Item(27).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>y</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:e><m:sup><m:r><m:t>5</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath29 = $d.OMaths.Item(29).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>7</m:t></m:r><m:sSup><m:e><m:r><m:t>z</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath31 = $d.OMaths.Item(31).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>5</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath35 = $d.OMaths.Item(35).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>2</m:t></m:r></m:num><m:den><m:r><m:t>a</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>a</m:t></m:r></m:num><m:den><m:r><m:t>12</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath37 = $d.OMaths.Item(37).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>a</m:t></m:r></m:num><m:den><m:r><m:t>b</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>4</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>c</m:t></m:r></m:num><m:den><m:r><m:t>d</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>e</m:t></m:r></m:num><m:den><m:r><m:t>f</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup></m:oMath>')
$xOMath39 = $d.OMaths.Item(39).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>a</m:t></m:r></m:e><m:sup><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:sup></m:sSup></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>b</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:sup></m:sSup></m:e></m:d></m:oMath>')
$xOMath40 = $d.OMaths.Item(40).Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>a</m:t></m:r></m:num><m:den><m:r><m:t>b</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>n</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>⋅</m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>c</m:t></m:r></m:num><m:den><m:r><m:t>d</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>n</m:t></m:r></m:sup></m:sSup></m:oMath>')

Write-Host "edit complete"
